$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Delete the two "Note 2" paragraphs about border-radius / <table>
#    (and the blank paragraph between them). They sit between the
#    blank paragraph that follows "Note : ... overflow: hidden." and
#    the final "Note 3 : ..." paragraph.
# --------------------------------------------------------------------
$pBorderRadius = $d.Paragraphs.Item(9)
$pTable        = $d.Paragraphs.Item(11)
$killRange = $d.Range($pBorderRadius.Range.Start, $pTable.Range.End)
$killRange.Delete()

# --------------------------------------------------------------------
# 2) The blank paragraph that used to precede the deleted text was
#    sized 36/36 (18pt); it now needs to match the surrounding 44/44
#    (22pt) paragraph-mark formatting. Its Range has no runs, so we
#    temporarily insert a placeholder character, size it, then remove
#    the character again - this leaves the paragraph-mark run
#    properties (w:pPr/w:rPr) updated without adding visible text.
# --------------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(8)
$blankStart = $pBlank.Range.Start
$placeholder = $d.Range($blankStart, $blankStart)
$placeholder.InsertAfter("X")
$pBlank = $d.Paragraphs.Item(8)
$pBlank.Range.Font.Size = 22
$charRange = $d.Range($pBlank.Range.Start, $pBlank.Range.Start + 1)
$charRange.Delete()

# --------------------------------------------------------------------
# 3) The final paragraph used to read "Note 3 : The overflow: ...".
#    It becomes "Note 2 : The overflow: ..." - note the run split
#    changes too: "Note 2" is its own run, immediately followed by
#    the relocated "_GoBack" bookmark, then a run starting " : The...".
#    The lastRenderedPageBreak and the proofErr wrapping "3 :" go away
#    as part of this text being retyped.
# --------------------------------------------------------------------
$pFinal = $d.Paragraphs.Item(9)
$finalStart = $pFinal.Range.Start
$headRange = $d.Range($finalStart, $finalStart + 9)   # "Note 3 : "
$headRange.Text = "Note 2 : "

# Move the _GoBack bookmark so it sits right after "Note 2"
$bookmarkPos = $finalStart + 6
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

Write-Output "done"
